# Jogos_da_Semana_FlashScore_2024-10-23.xlsx - odds refresh + new fixture
#
# 1) Rows 2-7: numeric odds columns (G:BD) get refreshed prices for the
#    already-listed fixtures (108 individual cell updates).
# 2) A new fixture (Slask Wroclaw vs Stal Mielec, id EXsrs05T) is inserted
#    as row 11, pushing the existing row 11 (tl3qMzJ5, FC Porto B vs Mafra)
#    down to row 12 unchanged. The sheet grows from A1:BD11 to A1:BD12.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Update odds values in rows 2-7 (existing matches) ---
# Row 2
$ws.Cells.Item(2, 17).Value = 2.2   # Q2
$ws.Cells.Item(2, 18).Value = 1.62   # R2
# Row 3
$ws.Cells.Item(3, 13).Value = 1.02   # M3
$ws.Cells.Item(3, 14).Value = 13.9   # N3
# Row 4
$ws.Cells.Item(4, 8).Value = 3.2   # H4
$ws.Cells.Item(4, 10).Value = 3.45   # J4
$ws.Cells.Item(4, 11).Value = 2.15   # K4
$ws.Cells.Item(4, 12).Value = 2.65   # L4
$ws.Cells.Item(4, 23).Value = 11.75   # W4
$ws.Cells.Item(4, 24).Value = 17.5   # X4
$ws.Cells.Item(4, 27).Value = 19   # AA4
$ws.Cells.Item(4, 28).Value = 18.5   # AB4
$ws.Cells.Item(4, 29).Value = 13   # AC4
$ws.Cells.Item(4, 30).Value = 5.9   # AD4
$ws.Cells.Item(4, 32).Value = 26   # AF4
$ws.Cells.Item(4, 34).Value = 8.5   # AH4
$ws.Cells.Item(4, 35).Value = 10.5   # AI4
$ws.Cells.Item(4, 37).Value = 18   # AK4
$ws.Cells.Item(4, 38).Value = 12.5   # AL4
$ws.Cells.Item(4, 39).Value = 16   # AM4
$ws.Cells.Item(4, 41).Value = 16   # AO4
$ws.Cells.Item(4, 42).Value = 18   # AP4
$ws.Cells.Item(4, 43).Value = 70   # AQ4
$ws.Cells.Item(4, 44).Value = 75   # AR4
$ws.Cells.Item(4, 45).Value = 150   # AS4
$ws.Cells.Item(4, 48).Value = 37   # AV4
$ws.Cells.Item(4, 50).Value = 11   # AX4
$ws.Cells.Item(4, 51).Value = 15.5   # AY4
$ws.Cells.Item(4, 52).Value = 40   # AZ4
$ws.Cells.Item(4, 53).Value = 55   # BA4
$ws.Cells.Item(4, 54).Value = 150   # BB4
# Row 5
$ws.Cells.Item(5, 7).Value = 2.42   # G5
$ws.Cells.Item(5, 10).Value = 2.95   # J5
$ws.Cells.Item(5, 12).Value = 3.2   # L5
$ws.Cells.Item(5, 16).Value = 3.74   # P5
$ws.Cells.Item(5, 17).Value = 1.65   # Q5
$ws.Cells.Item(5, 18).Value = 1.98   # R5
$ws.Cells.Item(5, 21).Value = 1.56   # U5
$ws.Cells.Item(5, 22).Value = 2.36   # V5
$ws.Cells.Item(5, 23).Value = 8.5   # W5
$ws.Cells.Item(5, 24).Value = 11.75   # X5
$ws.Cells.Item(5, 27).Value = 15   # AA5
$ws.Cells.Item(5, 28).Value = 18   # AB5
$ws.Cells.Item(5, 29).Value = 11.25   # AC5
$ws.Cells.Item(5, 32).Value = 30   # AF5
$ws.Cells.Item(5, 34).Value = 9   # AH5
$ws.Cells.Item(5, 35).Value = 13   # AI5
$ws.Cells.Item(5, 37).Value = 27   # AK5
$ws.Cells.Item(5, 38).Value = 17   # AL5
$ws.Cells.Item(5, 39).Value = 19   # AM5
$ws.Cells.Item(5, 41).Value = 13   # AO5
$ws.Cells.Item(5, 42).Value = 17.5   # AP5
$ws.Cells.Item(5, 44).Value = 70   # AR5
$ws.Cells.Item(5, 46).Value = 2.87   # AT5
$ws.Cells.Item(5, 50).Value = 14.5   # AX5
$ws.Cells.Item(5, 51).Value = 19   # AY5
$ws.Cells.Item(5, 52).Value = 60   # AZ5
$ws.Cells.Item(5, 53).Value = 80   # BA5
# Row 6
$ws.Cells.Item(6, 7).Value = 3.25   # G6
$ws.Cells.Item(6, 9).Value = 2.4   # I6
$ws.Cells.Item(6, 10).Value = 3.75   # J6
$ws.Cells.Item(6, 11).Value = 2   # K6
$ws.Cells.Item(6, 12).Value = 3.2   # L6
$ws.Cells.Item(6, 13).Value = 1.08   # M6
$ws.Cells.Item(6, 14).Value = 7.5   # N6
$ws.Cells.Item(6, 15).Value = 1.4   # O6
$ws.Cells.Item(6, 16).Value = 2.75   # P6
$ws.Cells.Item(6, 17).Value = 2.25   # Q6
$ws.Cells.Item(6, 18).Value = 1.62   # R6
$ws.Cells.Item(6, 19).Value = 1.5   # S6
$ws.Cells.Item(6, 20).Value = 2.5   # T6
$ws.Cells.Item(6, 21).Value = 1.95   # U6
$ws.Cells.Item(6, 22).Value = 1.8   # V6
$ws.Cells.Item(6, 24).Value = 15   # X6
$ws.Cells.Item(6, 25).Value = 12   # Y6
$ws.Cells.Item(6, 26).Value = 34   # Z6
$ws.Cells.Item(6, 27).Value = 29   # AA6
$ws.Cells.Item(6, 29).Value = 7.5   # AC6
$ws.Cells.Item(6, 31).Value = 15   # AE6
$ws.Cells.Item(6, 33).Value = 351   # AG6
$ws.Cells.Item(6, 34).Value = 7   # AH6
$ws.Cells.Item(6, 35).Value = 11   # AI6
$ws.Cells.Item(6, 36).Value = 10   # AJ6
$ws.Cells.Item(6, 37).Value = 23   # AK6
$ws.Cells.Item(6, 41).Value = 19   # AO6
$ws.Cells.Item(6, 42).Value = 29   # AP6
$ws.Cells.Item(6, 43).Value = 51   # AQ6
$ws.Cells.Item(6, 44).Value = 81   # AR6
$ws.Cells.Item(6, 45).Value = 251   # AS6
$ws.Cells.Item(6, 46).Value = 2.5   # AT6
$ws.Cells.Item(6, 49).Value = 4.33   # AW6
$ws.Cells.Item(6, 50).Value = 15   # AX6
$ws.Cells.Item(6, 52).Value = 51   # AZ6
$ws.Cells.Item(6, 54).Value = 201   # BB6
# Row 7
$ws.Cells.Item(7, 8).Value = 3.4   # H7
$ws.Cells.Item(7, 14).Value = 12   # N7
$ws.Cells.Item(7, 17).Value = 1.8   # Q7
$ws.Cells.Item(7, 18).Value = 2   # R7
$ws.Cells.Item(7, 19).Value = 1.36   # S7
$ws.Cells.Item(7, 20).Value = 3   # T7
$ws.Cells.Item(7, 29).Value = 12   # AC7
$ws.Cells.Item(7, 34).Value = 9.5   # AH7
$ws.Cells.Item(7, 39).Value = 26   # AM7
$ws.Cells.Item(7, 44).Value = 67   # AR7
$ws.Cells.Item(7, 45).Value = 151   # AS7
$ws.Cells.Item(7, 46).Value = 3   # AT7
$ws.Cells.Item(7, 48).Value = 51   # AV7
$ws.Cells.Item(7, 55).Value = 501   # BC7

# --- 2a) Insert a new row at 11: shifts the existing row 11
#         (tl3qMzJ5 / FC Porto B vs Mafra) down to row 12 ---
$ws.Rows.Item(11).Insert()

# --- 2b) Write the new fixture into row 11 (EXsrs05T / Slask Wroclaw vs Stal Mielec) ---
$ws.Cells.Item(11, 1).Value = "EXsrs05T"   # A11
$ws.Cells.Item(11, 2).Value = "23/10/2024"   # B11
$ws.Cells.Item(11, 3).Value = "13:30"   # C11
$ws.Cells.Item(11, 4).Value = "POLAND - EKSTRAKLASA"   # D11
$ws.Cells.Item(11, 5).Value = "Slask Wroclaw"   # E11
$ws.Cells.Item(11, 6).Value = "Stal Mielec"   # F11
$ws.Cells.Item(11, 7).Value = 1.9   # G11
$ws.Cells.Item(11, 8).Value = 3.6   # H11
$ws.Cells.Item(11, 9).Value = 3.8   # I11
$ws.Cells.Item(11, 10).Value = 2.6   # J11
$ws.Cells.Item(11, 11).Value = 2.2   # K11
$ws.Cells.Item(11, 12).Value = 4.33   # L11
$ws.Cells.Item(11, 13).Value = 1.05   # M11
$ws.Cells.Item(11, 14).Value = 11   # N11
$ws.Cells.Item(11, 15).Value = 1.3   # O11
$ws.Cells.Item(11, 16).Value = 3.4   # P11
$ws.Cells.Item(11, 17).Value = 2.03   # Q11
$ws.Cells.Item(11, 18).Value = 1.83   # R11
$ws.Cells.Item(11, 19).Value = 1.4   # S11
$ws.Cells.Item(11, 20).Value = 2.75   # T11
$ws.Cells.Item(11, 21).Value = 1.8   # U11
$ws.Cells.Item(11, 22).Value = 1.95   # V11
$ws.Cells.Item(11, 23).Value = 7   # W11
$ws.Cells.Item(11, 24).Value = 9   # X11
$ws.Cells.Item(11, 25).Value = 8.5   # Y11
$ws.Cells.Item(11, 26).Value = 17   # Z11
$ws.Cells.Item(11, 27).Value = 15   # AA11
$ws.Cells.Item(11, 28).Value = 26   # AB11
$ws.Cells.Item(11, 29).Value = 10   # AC11
$ws.Cells.Item(11, 30).Value = 7   # AD11
$ws.Cells.Item(11, 31).Value = 15   # AE11
$ws.Cells.Item(11, 32).Value = 51   # AF11
$ws.Cells.Item(11, 33).Value = 251   # AG11
$ws.Cells.Item(11, 34).Value = 11   # AH11
$ws.Cells.Item(11, 35).Value = 19   # AI11
$ws.Cells.Item(11, 36).Value = 13   # AJ11
$ws.Cells.Item(11, 37).Value = 41   # AK11
$ws.Cells.Item(11, 38).Value = 34   # AL11
$ws.Cells.Item(11, 39).Value = 41   # AM11
$ws.Cells.Item(11, 40).Value = 4   # AN11
$ws.Cells.Item(11, 41).Value = 10   # AO11
$ws.Cells.Item(11, 42).Value = 21   # AP11
$ws.Cells.Item(11, 43).Value = 34   # AQ11
$ws.Cells.Item(11, 44).Value = 51   # AR11
$ws.Cells.Item(11, 45).Value = 151   # AS11
$ws.Cells.Item(11, 46).Value = 2.75   # AT11
$ws.Cells.Item(11, 47).Value = 8   # AU11
$ws.Cells.Item(11, 48).Value = 51   # AV11
$ws.Cells.Item(11, 49).Value = 6   # AW11
$ws.Cells.Item(11, 50).Value = 21   # AX11
$ws.Cells.Item(11, 51).Value = 29   # AY11
$ws.Cells.Item(11, 52).Value = 67   # AZ11
$ws.Cells.Item(11, 53).Value = 101   # BA11
$ws.Cells.Item(11, 54).Value = 201   # BB11
$ws.Cells.Item(11, 55).Value = 81   # BC11
$ws.Cells.Item(11, 56).Value = 81   # BD11
